$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value
$ws.Range("C2").Value = 1000000

# Row 3 - 1 thread
$ws.Range("B3").Value = "0.000648"
$ws.Range("C3").Value = "0.332334"

# Row 4 - 2 threads
$ws.Range("B4").Value = "0.000714"
$ws.Range("C4").Value = "0.207382"

# Row 5 - 5 threads
$ws.Range("B5").Value = "0.000832"
$ws.Range("C5").Value = "0.100345"

# Row 6 - 15 threads
$ws.Range("B6").Value = "0.002012"
$ws.Range("C6").Value = "0.090462"

# Row 7 - 25 threads
$ws.Range("B7").Value = "0.003812"
$ws.Range("C7").Value = "0.107350"

# Row 8 - single thread (new row)
$ws.Range("A8").Value = "single thread"
$ws.Range("B8").Value = "0.000205"
$ws.Range("C8").Value = "0.185292 "

# Selection moved to B7
$ws.Range("B7").Select()
